# Adds season-record columns (Wins / Losses / Ties) to the roster sheet.
# New columns: AD = Wins, AE = Losses, AF = Ties.
# Header row (row 1) gets the same bold/centered/bordered style as the
# existing "Unnamed: 28" header in AC1; data rows 2-54 get the team's
# season record (51 wins, 111 losses, 0 ties) as plain numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
# Copy AC1's formatting (bold font, thin border, center/top alignment)
# onto the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows -----------------------------------------------------------
# Every player row shares the same team season record.
$wins = 51
$losses = 111
$ties = 0

for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
